$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K data shifts to F:M),
# matching the quarterly-report convention of adding the two most
# recent quarters as the new leftmost data columns.
$ws.Range("D7:E102").Insert(-4161)

# Copy number formatting from the (now-shifted) old column D, currently
# column F, into the two newly inserted columns so the new cells keep
# the same style as the rest of each row (date format row 7/38/80,
# thousands format elsewhere).
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1310100
$ws.Range("E8").Value = 1283800
$ws.Range("D9").Value = 670900
$ws.Range("E9").Value = 660300
$ws.Range("D10").Value = 639200
$ws.Range("E10").Value = 623500
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 12600
$ws.Range("E14").Value = -2600
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 1051500
$ws.Range("E17").Value = 1016500
$ws.Range("D18").Value = 258600
$ws.Range("E18").Value = 267300
$ws.Range("D20").Value = 7500
$ws.Range("E20").Value = 6600
$ws.Range("D21").Value = 572900
$ws.Range("E21").Value = 579600
$ws.Range("D22").Value = 130000
$ws.Range("E22").Value = 130600
$ws.Range("D23").Value = 136100
$ws.Range("E23").Value = 143300
$ws.Range("D24").Value = 32600
$ws.Range("E24").Value = 18500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 103500
$ws.Range("E26").Value = 124800
$ws.Range("D27").Value = 103500
$ws.Range("E27").Value = 124800
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 6500
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -7500
$ws.Range("E32").Value = -6600
$ws.Range("D33").Value = 110000
$ws.Range("E33").Value = 124800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 110000
$ws.Range("E35").Value = 124800
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 606200
$ws.Range("E41").Value = 870500
$ws.Range("D42").Value = 4500
$ws.Range("E42").Value = 15400
$ws.Range("D43").Value = 630100
$ws.Range("E43").Value = 662400
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 274900
$ws.Range("E45").Value = 258700
$ws.Range("D46").Value = 1515700
$ws.Range("E46").Value = 1807000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 11026000
$ws.Range("E48").Value = 10682800
$ws.Range("D49").Value = 7169700
$ws.Range("E49").Value = 7235900
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 533300
$ws.Range("E52").Value = 562300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 20244600
$ws.Range("E54").Value = 20288100
$ws.Range("D57").Value = 756700
$ws.Range("E57").Value = 739100
$ws.Range("D58").Value = 452000
$ws.Range("E58").Value = 322100
$ws.Range("D59").Value = 306400
$ws.Range("E59").Value = 400100
$ws.Range("D60").Value = 1515100
$ws.Range("E60").Value = 1461300
$ws.Range("D61").Value = 10880500
$ws.Range("E61").Value = 11032500
$ws.Range("D62").Value = 629800
$ws.Range("E62").Value = 634100
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 13025400
$ws.Range("E66").Value = 13127900
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 889900
$ws.Range("E72").Value = 779900
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 7219300
$ws.Range("E76").Value = 7160200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 110000
$ws.Range("E81").Value = 124800
$ws.Range("D83").Value = 306900
$ws.Range("E83").Value = 305700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 558600
$ws.Range("E89").Value = 417200
$ws.Range("D91").Value = -680700
$ws.Range("E91").Value = -545500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -725600
$ws.Range("E94").Value = -620000
$ws.Range("D96").Value = -183900
$ws.Range("E96").Value = -186000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -100900
$ws.Range("E100").Value = 104100
$ws.Range("D101").Value = -3000
$ws.Range("E101").Value = -5100
$ws.Range("D102").Value = -270800
$ws.Range("E102").Value = -103800
